# Update column F (dSF) values for specific rows per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 6
$ws.Range("F4").Value = -1
$ws.Range("F10").Value = -1
$ws.Range("F16").Value = -4
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 1
$ws.Range("F23").Value = -1
